$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): extend with P1=14 and Q1=15, copying the existing
# header style (bold font, thin border, centered) from O1 so the new
# cells pick up style index 1 just like the rest of the header.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Data rows 2-25: the I/K and M/O columns swap their 1/2 values, and two
# new unstyled columns P and Q (value 2) are appended to every row.
for ($r = 2; $r -le 25; $r++) {
    $ws.Range("I$r").Value = 2
    $ws.Range("K$r").Value = 1
    $ws.Range("M$r").Value = 2
    $ws.Range("O$r").Value = 1
    $ws.Range("P$r").Value = 2
    $ws.Range("Q$r").Value = 2
}
